# Applies the numeric updates to Sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# as produced by the scheduled Phoenix_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 699.5454999999999
$ws.Range("I19").Value = 499.5
$ws.Range("J19").Value = 939.6
$ws.Range("K19").Value = 499.5
$ws.Range("L19").Value = 939.6
$ws.Range("M19").Value = -324.5
$ws.Range("N19").Value = -1289.6

$ws.Range("H28").Value = 915.125
$ws.Range("I28").Value = 1086.4166
$ws.Range("J28").Value = 401.25
$ws.Range("K28").Value = 1086.4166
$ws.Range("L28").Value = 401.25
$ws.Range("M28").Value = -601.4166
$ws.Range("N28").Value = -1371.25

$ws.Range("H113").Value = 2620.3333
$ws.Range("I113").Value = 2439
$ws.Range("J113").Value = 2892.3333
$ws.Range("K113").Value = 2439
$ws.Range("L113").Value = 2892.3333
$ws.Range("M113").Value = 815
$ws.Range("N113").Value = -9400.3333

$ws.Range("H121").Value = 2690.5
$ws.Range("J121").Value = 2690.5
$ws.Range("L121").Value = 8071.5
$ws.Range("N121").Value = -11565.5

$ws.Range("H132").Value = 2103.254
$ws.Range("I132").Value = 2121.6047
$ws.Range("J132").Value = 2063.8
$ws.Range("K132").Value = 6364.8141
$ws.Range("L132").Value = 6191.400000000001
$ws.Range("M132").Value = -3834.8141
$ws.Range("N132").Value = -11251.4

$ws.Range("H138").Value = 2657.25
$ws.Range("I138").Value = 900
$ws.Range("J138").Value = 3008.7
$ws.Range("K138").Value = 2700
$ws.Range("L138").Value = 9026.099999999999
$ws.Range("M138").Value = 2440
$ws.Range("N138").Value = -19306.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2995.3098
$ws.Range("I32").Value = 2817.2207
$ws.Range("J32").Value = 7032
$ws.Range("K32").Value = 2817.2207
$ws.Range("L32").Value = 7032
$ws.Range("M32").Value = -2530.2207
$ws.Range("N32").Value = -7606

$ws.Range("H61").Value = 2763.3572
$ws.Range("I61").Value = 2333.6572
$ws.Range("J61").Value = 4911.857
$ws.Range("K61").Value = 2333.6572
$ws.Range("L61").Value = 4911.857
$ws.Range("M61").Value = -2121.6572
$ws.Range("N61").Value = -5335.857

$ws.Range("H63").Value = 1599
$ws.Range("I63").Value = 1599
$ws.Range("K63").Value = 1599
$ws.Range("M63").Value = -913

$ws.Range("H66").Value = 1599
$ws.Range("I66").Value = 1599
$ws.Range("K66").Value = 7995
$ws.Range("M66").Value = -4563

$ws.Range("H74").Value = 1929.5
$ws.Range("I74").Value = 1759.2307
$ws.Range("J74").Value = 2482.875
$ws.Range("K74").Value = 1759.2307
$ws.Range("L74").Value = 2482.875
$ws.Range("M74").Value = -885.2307000000001
$ws.Range("N74").Value = -4230.875

$ws.Range("H77").Value = 1929.5
$ws.Range("I77").Value = 1759.2307
$ws.Range("J77").Value = 2482.875
$ws.Range("K77").Value = 8796.1535
$ws.Range("L77").Value = 12414.375
$ws.Range("M77").Value = -4428.1535
$ws.Range("N77").Value = -21150.375

$ws.Range("H136").Value = 2763.3572
$ws.Range("I136").Value = 2333.6572
$ws.Range("J136").Value = 4911.857
$ws.Range("K136").Value = 7000.971600000001
$ws.Range("L136").Value = 14735.571
$ws.Range("M136").Value = -4450.971600000001
$ws.Range("N136").Value = -19835.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2283.1538
$ws.Range("I94").Value = 2334.8
$ws.Range("J94").Value = 2111
$ws.Range("K94").Value = 2334.8
$ws.Range("L94").Value = 2111
$ws.Range("M94").Value = -1883.8
$ws.Range("N94").Value = -3013

$ws.Range("H99").Value = 1912.9286
$ws.Range("I99").Value = 2032
$ws.Range("J99").Value = 1198.5
$ws.Range("K99").Value = 2032
$ws.Range("L99").Value = 1198.5
$ws.Range("M99").Value = -534
$ws.Range("N99").Value = -4194.5

$ws.Range("H106").Value = 46999.5
$ws.Range("J106").Value = 46999.5
$ws.Range("L106").Value = 46999.5
$ws.Range("N106").Value = -49523.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1666.7778
$ws.Range("I132").Value = 1643.6666
$ws.Range("J132").Value = 1851.6666
$ws.Range("K132").Value = 4930.9998
$ws.Range("L132").Value = 5554.9998
$ws.Range("M132").Value = -2400.9998
$ws.Range("N132").Value = -10614.9998

$ws.Range("H134").Value = 2012.5714
$ws.Range("I134").Value = 2014.579
$ws.Range("J134").Value = 1993.5
$ws.Range("K134").Value = 6043.737
$ws.Range("L134").Value = 5980.5
$ws.Range("M134").Value = -3508.737
$ws.Range("N134").Value = -11050.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 83341680
$ws.Range("I55").Value = 724.6667
$ws.Range("J55").Value = 111122000
$ws.Range("K55").Value = 2174.0001
$ws.Range("L55").Value = 333366000
$ws.Range("M55").Value = -1997.0001
$ws.Range("N55").Value = -333366354

$ws.Range("H62").Value = 8344.130999999999
$ws.Range("I62").Value = 5333
$ws.Range("J62").Value = 8795.799999999999
$ws.Range("K62").Value = 15999
$ws.Range("L62").Value = 26387.4
$ws.Range("M62").Value = -15313
$ws.Range("N62").Value = -27759.4

$ws.Range("H65").Value = 8344.130999999999
$ws.Range("I65").Value = 5333
$ws.Range("J65").Value = 8795.799999999999
$ws.Range("K65").Value = 47997
$ws.Range("L65").Value = 79162.2
$ws.Range("M65").Value = -44565
$ws.Range("N65").Value = -86026.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4966.628
$ws.Range("I102").Value = 4746.879
$ws.Range("J102").Value = 5691.8
$ws.Range("K102").Value = 4746.879
$ws.Range("L102").Value = 5691.8
$ws.Range("M102").Value = -3124.879
$ws.Range("N102").Value = -8935.799999999999

$ws.Range("H109").Value = 39234.5
$ws.Range("J109").Value = 39234.5
$ws.Range("L109").Value = 39234.5
$ws.Range("N109").Value = -41314.5

$ws.Range("H113").Value = 5770.5864
$ws.Range("I113").Value = 6654.304
$ws.Range("J113").Value = 2383
$ws.Range("K113").Value = 6654.304
$ws.Range("L113").Value = 2383
$ws.Range("M113").Value = -4484.304
$ws.Range("N113").Value = -6723

$ws.Range("H122").Value = 60146.11
$ws.Range("I122").Value = 75768.36
$ws.Range("J122").Value = 5468.25
$ws.Range("K122").Value = 227305.08
$ws.Range("L122").Value = 16404.75
$ws.Range("M122").Value = -224855.08
$ws.Range("N122").Value = -21304.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4874.625
$ws.Range("I7").Value = 4143.143
$ws.Range("J7").Value = 9995
$ws.Range("K7").Value = 4143.143
$ws.Range("L7").Value = 9995
$ws.Range("M7").Value = -4031.143
$ws.Range("N7").Value = -10219

$ws.Range("H25").Value = 25503.5
$ws.Range("I25").Value = 25503.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 25503.5
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -25273.5
$ws.Range("N25").ClearContents()

$ws.Range("H122").Value = 3083.6
$ws.Range("I122").Value = 3094.027
$ws.Range("J122").Value = 2955
$ws.Range("K122").Value = 9282.081
$ws.Range("L122").Value = 8865
$ws.Range("M122").Value = -6832.081
$ws.Range("N122").Value = -13765

$ws.Range("H126").Value = 4874.625
$ws.Range("I126").Value = 4143.143
$ws.Range("J126").Value = 9995
$ws.Range("K126").Value = 12429.429
$ws.Range("L126").Value = 29985
$ws.Range("M126").Value = -9959.429
$ws.Range("N126").Value = -34925

$ws.Range("H140").Value = 83038.39999999999
$ws.Range("J140").Value = 83038.39999999999
$ws.Range("L140").Value = 83038.39999999999
$ws.Range("N140").Value = -93398.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 23585.375
$ws.Range("J41").Value = 23666.5
$ws.Range("L41").Value = 23666.5
$ws.Range("N41").Value = -24446.5

$ws.Range("H52").Value = 38018.152
$ws.Range("I52").Value = 20099.334
$ws.Range("J52").Value = 40355.39
$ws.Range("K52").Value = 20099.334
$ws.Range("L52").Value = 40355.39
$ws.Range("M52").Value = -19873.334
$ws.Range("N52").Value = -40807.39

$ws.Range("H136").Value = 5129012.5
$ws.Range("I136").Value = 5917880
$ws.Range("J136").Value = 1372.5
$ws.Range("K136").Value = 17753640
$ws.Range("L136").Value = 4117.5
$ws.Range("M136").Value = -17751090
$ws.Range("N136").Value = -9217.5
